$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new "season record" header cells right after the
# existing last column (AC) -> AD, AE, AF.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new header cells the same formatting as the rest of the
# header row (bold, centered, bordered) by copying the style from the
# neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Populate every data row (2-49) with the team's season record: the
# same Wins/Losses/Ties values apply to every player on the roster.
$ws.Range("AD2:AD49").Value = 78
$ws.Range("AE2:AE49").Value = 84
$ws.Range("AF2:AF49").Value = 0
